$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill the "Approved/Rejected" column (I) with "Approved" for most rows,
# but row 17 is "Reject" with a reason ("wrong data") in column J.
$approvedRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,18,19,20,21)
foreach ($r in $approvedRows) {
    $ws.Cells.Item($r, 9).Value = "Approved"
}

$ws.Cells.Item(17, 9).Value = "Reject"
$ws.Cells.Item(17, 10).Value = "wrong data"

# Update the view: scrolled/selected cell moved from D1/J22 to F1/J17.
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("J17").Select()
